# LKS.xlsx player-stats update (2022-10-01 snapshot).
# The sheet stores every stat as text (t="inlineStr" in the original OOXML),
# so numeric-looking values are written with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# Row 2 - D. Arndt
Set-TextValue "E2" "856"
Set-TextValue "F2" "10"
Set-TextValue "G2" "10"
Set-TextValue "I2" "1"

# Row 3 - A. Bobek
Set-TextValue "E3" "134"
Set-TextValue "F3" "2"
Set-TextValue "H3" "1"
Set-TextValue "J3" "5"

# Row 7 - M. Dabrowski
Set-TextValue "E7" "990"
Set-TextValue "F7" "11"
Set-TextValue "G7" "11"
Set-TextValue "L7" "6"

# Row 8 - K. Dankowski
Set-TextValue "E8" "833"
Set-TextValue "F8" "12"
Set-TextValue "G8" "8"
Set-TextValue "L8" "1"

# Row 9 - O. Koprowski
Set-TextValue "E9" "364"
Set-TextValue "F9" "6"
Set-TextValue "G9" "4"

# Row 10 - M. Lorenc
Set-TextValue "E10" "507"
Set-TextValue "F10" "8"
Set-TextValue "G10" "6"

# Row 11 - A. Marciniak
Set-TextValue "J11" "4"

# Row 12 - Nacho Monsalve
Set-TextValue "E12" "1080"
Set-TextValue "F12" "12"
Set-TextValue "G12" "12"

# Row 14 - A. Tutyskinas
Set-TextValue "J14" "3"

# Row 15 - M. Wszolek
Set-TextValue "J15" "8"

# Row 16 - B. Biel
Set-TextValue "E16" "636"
Set-TextValue "F16" "12"
Set-TextValue "G16" "8"
Set-TextValue "I16" "7"

# Row 18 - D. Kort
Set-TextValue "E18" "697"
Set-TextValue "F18" "11"
Set-TextValue "G18" "9"
Set-TextValue "K18" "1"
Set-TextValue "L18" "3"

# Row 20 - J. Kuzma
Set-TextValue "E20" "326"
Set-TextValue "F20" "8"
Set-TextValue "G20" "3"
Set-TextValue "I20" "3"

# Row 22 - D. Nowacki
Set-TextValue "J22" "3"

# Row 23 - V. Okhronchuk
Set-TextValue "J23" "10"

# Row 26 - M. Trabka
Set-TextValue "E26" "956"
Set-TextValue "F26" "12"
Set-TextValue "G26" "11"

# Row 27 - N. Balongo
Set-TextValue "E27" "706"
Set-TextValue "F27" "12"
Set-TextValue "H27" "3"
Set-TextValue "J27" "3"

# Row 28 - G. Glapka
Set-TextValue "E28" "12"
Set-TextValue "F28" "2"
Set-TextValue "H28" "2"
Set-TextValue "J28" "2"

# Row 29 - P. Janczukowicz
Set-TextValue "E29" "312"
Set-TextValue "F29" "10"
Set-TextValue "H29" "8"
Set-TextValue "J29" "9"
Set-TextValue "L29" "2"

# Row 31 - M. Radaszkiewicz
Set-TextValue "E31" "100"
Set-TextValue "F31" "2"
Set-TextValue "G31" "1"
Set-TextValue "I31" "1"
